$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "C" (date header / rating
# change) column. This shifts the old C column -> E; column B is untouched
# by the insert itself.
$ws.Columns("C:D").Insert()

# Re-assert the column widths (8 "characters", same as the pre-existing custom
# width used on the old date column) so C, D and E all keep an explicit
# custom width of 8 after the shift.
$ws.Columns("C").ColumnWidth = 7.14
$ws.Columns("D").ColumnWidth = 7.14
$ws.Columns("E").ColumnWidth = 7.14

# New header row: two new date columns (Jun_15, Jun_17) plus re-labeling of
# the existing header cells. Order matters for shared-string placement: C1
# ("Jun_15") is written before B1 ("Jun_17") so the new unique strings land
# in the same order as the target workbook.
$ws.Cells.Item(1,3).Value = "Jun_15"
$ws.Cells.Item(1,4).Value = "Jun_13"
$ws.Cells.Item(1,2).Value = "Jun_17"

# Fill the two new body columns (C, D) with the same "UN" placeholder value
# already used in column B, for every data row.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r,3).Value = "UN"
    $ws.Cells.Item($r,4).Value = "UN"
}
